$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Add 1 more working day ("cong") for Nguyen Huu Quang, flowing through the
# derived salary totals.
$ws.Range("B1").Value = 18
$ws.Range("B12").Value = 27
$ws.Range("B13").Value = 3857142.857142857
$ws.Range("B32").Value = 1557142.857142857
$ws.Range("B34").Value = 1557142.857142857
